$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H11/I11: status changed from "Falta" to "Completado", date 44153 -> 44155
# Copy the formatting already used for "Completado" cells (e.g. H3) so the
# shared style is reused instead of creating a new cell style.
$ws.Range("H3").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Value = "Completado"
$ws.Range("I11").Value = 44155

# Update H12/I12: same change
$ws.Range("H3").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H12").Value = "Completado"
$ws.Range("I12").Value = 44155

# Update I15, I35-I39 dates only
$ws.Range("I15").Value = 44155
$ws.Range("I35").Value = 44155
$ws.Range("I36").Value = 44155
$ws.Range("I37").Value = 44155
$ws.Range("I38").Value = 44155
$ws.Range("I39").Value = 44155

# Update sheet view: remove topLeftCell freeze position, change selection to I15
$ws.Range("I15").Select()
